$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.931.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = "'2.877.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.40%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'587.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'138.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.13%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -3.34%  '
$ws.Range("D9").Value = "'6.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("D10").Value = "'0.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.71%  '
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("E12").Value = '  -3.87%  '
$ws.Range("D13").Value = "'32.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.04%  '
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").Value = "'3.348.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = "'60.882.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = "'2.865.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("E18").Value = '  -3.16%  '
$ws.Range("D19").Value = "'424.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").Value = "'13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").Value = "'0.652"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.88%  '
$ws.Range("D22").Value = "'6.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.39%  '
$ws.Range("D23").Value = "'79.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.27%  '
$ws.Range("D24").Value = "'10.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.22%  '
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'2.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.63%  '
$ws.Range("D27").Value = "'11.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.46%  '
$ws.Range("E28").Value = '  -3.18%  '
$ws.Range("E29").Value = '  -9.22%  '
$ws.Range("D30").Value = "'6.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.68%  '
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").Value = "'25.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.31%  '
$ws.Range("E33").Value = '  -4.23%  '
$ws.Range("D34").Value = "'0.0₃0846"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("D35").Value = "'0.964"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.49%  '
$ws.Range("E36").Value = '  -3.85%  '
$ws.Range("D37").Value = "'48.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("D38").Value = "'2.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.11%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = "'1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.39%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("D42").Value = "'38.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.00%  '
$ws.Range("D43").Value = "'0.262"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.19%  '
$ws.Range("D44").Value = "'2.655.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").Value = "'132.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.60%  '
$ws.Range("E46").Value = '  -3.86%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value = "'342.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.79%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").Value = "'22.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.99%  '
$ws.Range("D50").Value = "'0.102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.00%  '
$ws.Range("D51").Value = "'1.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.23%  '
